# 1 add testcase(face_detect) 2 enhance clus_face update db
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: add a new "3#" test-machine row (row 12) ---
$ws1.Columns.Item(3).ColumnWidth = 25.28

$ws1.Range("B11").Copy()
$ws1.Range("B12").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("C11").Copy()
$ws1.Range("D12").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("B12").Value = "3#"
$ws1.Range("D12").Value = 23.71

# --- Sheet2: build out the face-detect testcase table ---
$ws1.Range("B2").Copy()
$ws2.Range("A2:J16").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A2:J16").RowHeight = 15.75

$ws2.Columns.Item(2).ColumnWidth = 11.57
$ws2.Columns.Item(3).ColumnWidth = 14.85
$ws2.Columns.Item(5).ColumnWidth = 10.7

# Row 2 - header
$ws2.Range("B2").Value = "No."
$ws2.Range("C2").Value = "width"
$ws2.Range("D2").Value = "height"
$ws2.Range("E2").Value = "detect sec"
$ws2.Range("F2").Value = "file size"

# Row 3
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 7952
$ws2.Range("D3").Value = 5304
$ws2.Range("E3").Value = 37
$ws2.Range("F3").Value = "1MB"

# Row 4
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = 8640
$ws2.Range("D4").Value = 5760
$ws2.Range("E4").Value = 43
$ws2.Range("F4").Value = "16MB"

# Row 5
$ws2.Range("B5").Value = 3
$ws2.Range("C5").Value = 1920
$ws2.Range("D5").Value = 1080
$ws2.Range("E5").Value = 1.8
$ws2.Range("F5").Value = "100k-300kB"

# Row 6
$ws2.Range("B6").Value = 4
$ws2.Range("C6").Value = 1280
$ws2.Range("D6").Value = 720
$ws2.Range("E6").Value = 0.9
$ws2.Range("F6").Value = "60k-300k"

# Row 7
$ws2.Range("B7").Value = 5
$ws2.Range("C7").Value = 720
$ws2.Range("D7").Value = 480
$ws2.Range("E7").Value = 0.4
$ws2.Range("F7").Value = "20k-90k"

# Row 9
$ws2.Range("C9").Value = "129 pics football"

# Row 11 - test machine info
$ws2.Range("B11").Value = "test machine"
$ws2.Range("C11").Value = "192.168.23.71"
$ws2.Range("D11").Value = "cpu"
$ws2.Range("E11").Value = "Intel(R) Xeon(R) CPU E5-2620 v2 @ 2.10GHz"
$ws2.Range("E11").HorizontalAlignment = -4131  # xlHAlignLeft
$ws2.Range("I11").Value = "24core"

# Row 12
$ws2.Range("D12").Value = "mem"
$ws2.Range("E12").Value = "32GB"

# --- selections / active sheet ---
[void]$ws1.Range("B2:F12").Select()
[void]$ws2.Range("F16").Select()
$ws2.Activate()

Write-Host "edit complete"
